$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Style templates come from rows 1-9 / row1 header, which are unchanged by this edit ----
# A-column style (bold, valign top)          -> template cell A3
# B-column style (valign top, wrap)            -> template cell B3
# C-column style (valign top, wrap, red font)  -> template cell C3
# Row-1 header style (bold, valign top, no wrap) -> template cell B1

# Clear the region that needs restructuring (rows 10-26); rows 1-9 stay untouched.
$ws.Range("A10:C26").Clear()

# Row 10
$ws.Range("A3").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B3").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = "Aplicação dos conhecimentos adquiridos nas disciplinas obrigatórias e das competências desenvolvidas durante o Curso a uma situação possível do ambiente profissional"
$ws.Range("C3").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = "Aplicação dos conhecimentos adquiridos nas disciplinas obrigatórias e das competências desenvolvidas durante o Curso a uma situação possível do ambiente profissional"
$ws.Rows.Item(10).AutoFit()
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Range("A3").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Objectives:"
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Range("A3").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Rows.Item(12).AutoFit()

# Row 13
$ws.Range("B3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C3").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Rows.Item(13).AutoFit()

# Row 14
$ws.Range("B3").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C3").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Rows.Item(14).AutoFit()

# Row 15
$ws.Range("B3").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C3").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Rows.Item(15).AutoFit()

# Row 16
$ws.Range("A3").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Programa resumido:"
$ws.Range("B3").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = "Elaboração de um projeto de engenharia, de pesquisa científica ou modelo de negócio, Desenvolvimento do projeto, com características inter e transdisciplinar"
$ws.Range("C3").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = "Elaboração de um projeto de engenharia, de pesquisa científica ou modelo de negócio, Desenvolvimento do projeto, com características inter e transdisciplinar"
$ws.Rows.Item(16).AutoFit()
$ws.Rows.Item(16).RowHeight = 60

# Row 17
$ws.Range("A3").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Short syllabus:"
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(17).RowHeight = 60

# Row 18
$ws.Range("A3").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Programa:"
$ws.Range("B3").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "A disciplina consiste no desenvolvimento de um projeto supervisionado por docente e/ou profissional de engenharia que poderá ser realizado em grupo ou de forma individual.1. Elaboração de um projeto de engenharia, ou pesquisa científica ou modelo de negócio utilizando as competências adquiridas nos Projetos de Engenharia I, II, III e IV). O projeto deve atender os princípios de planejamento e gestão de projetos ou de negócios, inclusive, se possível construindo modelo ou protótipo físico e/ou digital2. Desenvolvimento do Projeto – em projetos experimentais deverão ser produzidos alguns produtos, processos ou sistemas reais, teste de modelos ou protótipos.3. O aluno, individualmente ou em equipe, deverá elaborar uma monografia ou plano de negócio do projeto e submete-la a apreciação de uma banca – a monografia ou plano de negócio deve atender aos padrões estabelecidos e utilizados nas disciplinas de Projetos de Engenharia I, II, III e IV).4. Apresentação do Projeto Final de Curso para uma banca de três professores, sendo um orientador, no caso de ter mais de um, e dois outros membros, entre eles preferencialmente, um vindo da indústria do ramo de atividade do tema escolhido."
$ws.Range("C3").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "A disciplina consiste no desenvolvimento de um projeto supervisionado por docente e/ou profissional de engenharia que poderá ser realizado em grupo ou de forma individual.1. Elaboração de um projeto de engenharia, ou pesquisa científica ou modelo de negócio utilizando as competências adquiridas nos Projetos de Engenharia I, II, III e IV). O projeto deve atender os princípios de planejamento e gestão de projetos ou de negócios, inclusive, se possível construindo modelo ou protótipo físico e/ou digital2. Desenvolvimento do Projeto – em projetos experimentais deverão ser produzidos alguns produtos, processos ou sistemas reais, teste de modelos ou protótipos.3. O aluno, individualmente ou em equipe, deverá elaborar uma monografia ou plano de negócio do projeto e submete-la a apreciação de uma banca – a monografia ou plano de negócio deve atender aos padrões estabelecidos e utilizados nas disciplinas de Projetos de Engenharia I, II, III e IV).4. Apresentação do Projeto Final de Curso para uma banca de três professores, sendo um orientador, no caso de ter mais de um, e dois outros membros, entre eles preferencialmente, um vindo da indústria do ramo de atividade do tema escolhido."
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(18).RowHeight = 120

# Row 19
$ws.Range("A3").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Syllabus:"
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(19).RowHeight = 120

# Row 20
$ws.Range("A3").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Avaliação:"
$ws.Rows.Item(20).AutoFit()

# Row 21
$ws.Range("A3").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Método:"
$ws.Range("B3").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").Value = "O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas."
$ws.Range("C3").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C21").Value = "O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas."
$ws.Rows.Item(21).AutoFit()
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A3").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Critério:"
$ws.Range("B3").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = "A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina."
$ws.Range("C3").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina."
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$ws.Range("A3").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Norma de recuperação:"
$ws.Range("B3").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "não há"
$ws.Range("C3").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "não há"
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(23).RowHeight = 60

# Row 24
$ws.Range("A3").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").Value = "Bibliografia:"
$ws.Range("B3").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").Value = "A ser definida em função do projeto"
$ws.Range("C3").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = "A ser definida em função do projeto"
$ws.Rows.Item(24).AutoFit()
$ws.Rows.Item(24).RowHeight = 120

# Row 25
$ws.Range("A3").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = "Requisitos:"
$ws.Rows.Item(25).AutoFit()

# Row 26
$ws.Range("B3").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("B26").Value = "LOM3110 -  Projeto Integrado em Engenharia de Materiais III  (Requisito fraco)`n"
$ws.Range("C3").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = "LOM3110 -  Projeto Integrado em Engenharia de Materiais III  (Requisito fraco)`n"
$ws.Rows.Item(26).AutoFit()
$ws.Rows.Item(26).RowHeight = 30

$excel.CutCopyMode = $false